# Auto-generated Excel COM-interop edit script
# Applies the cell-value updates for the 'cryptos' worksheet refresh
# described by the commit 'Updated cryptos list on Sat Aug 24 22:49:26 UTC 2024 with GitHub Actions'.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text can be written directly (Excel will not mis-parse them as numbers).
$plainUpdates = @{
    'D2' = '63.894.11'
    'E2' = '  -0.38%  '
    'D3' = '2.749.76'
    'E3' = '  +0.18%  '
    'E5' = '  -2.48%  '
    'E6' = '  +4.24%  '
    'E7' = '  +0.39%  '
    'E8' = '  -0.20%  '
    'E9' = '  -1.56%  '
    'E10' = '  +0.65%  '
    'E11' = '  -15.77%  '
    'E12' = '  +0.17%  '
    'D13' = '3.235.75'
    'E13' = '  -0.06%  '
    'E14' = '  +1.16%  '
    'D15' = '63.800.34'
    'E15' = '  -0.35%  '
    'E16' = '  +1.41%  '
    'D17' = '2.747.82'
    'E17' = '  -0.89%  '
    'E18' = '  +1.64%  '
    'E19' = '  +0.78%  '
    'E20' = '  -1.55%  '
    'E21' = '  -2.20%  '
    'E22' = '  +5.80%  '
    'E23' = '  +0.49%  '
    'E24' = '  +0.06%  '
    'E25' = '  +2.16%  '
    'E26' = '  +0.09%  '
    'E27' = '  +0.46%  '
    'D28' = '0.0₃0935'
    'E28' = '  +1.74%  '
    'E29' = '  -1.66%  '
    'E30' = '  -0.07%  '
    'E31' = '  +5.63%  '
    'E32' = '  -2.21%  '
    'E33' = '  +0.21%  '
    'B34' = 'EthereumClassic'
    'C34' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'E34' = '  -0.65%  '
    'B35' = 'NEARProtocol'
    'C35' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'E35' = '  +4.40%  '
    'E36' = '  +2.78%  '
    'E37' = '  +0.79%  '
    'E38' = '  -1.06%  '
    'E39' = '  -0.99%  '
    'E40' = '  +8.54%  '
    'E41' = '  -5.04%  '
    'E42' = '  +0.61%  '
    'E43' = '  -0.90%  '
    'E44' = '  +1.12%  '
    'E45' = '  -0.59%  '
    'E46' = '  -0.02%  '
    'E47' = '  -1.43%  '
    'E48' = '  -4.30%  '
    'E49' = '  +0.55%  '
    'E50' = '  +0.53%  '
    'E51' = '  +0.85%  '
}

foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

# Cells whose new text is a plain decimal number (e.g. '1.00', '20.58').
# These must be forced to Text format first, otherwise Excel would silently
# reinterpret them as numeric values and mangle formatting/precision
# (e.g. '1.00' -> 1, '20.58' -> 20.579999999999998).
$textUpdates = @{
    'D4' = '1.00'
    'D5' = '579.46'
    'D6' = '159.09'
    'D8' = '0.611'
    'D10' = '0.392'
    'D11' = '5.70'
    'D12' = '0.159'
    'D14' = '27.03'
    'D18' = '12.29'
    'D21' = '6.89'
    'D22' = '0.569'
    'D23' = '0.998'
    'D24' = '66.37'
    'D26' = '8.69'
    'D27' = '1.00'
    'D30' = '7.14'
    'D31' = '1.26'
    'D32' = '168.50'
    'D34' = '20.58'
    'D35' = '5.02'
    'D36' = '1.48'
    'D40' = '6.17'
    'D41' = '333.05'
    'D42' = '39.72'
    'D43' = '22.06'
    'D44' = '0.0600'
    'D45' = '22.00'
    'D47' = '0.641'
    'D48' = '137.20'
    'D50' = '1.00'
    'D51' = '11.05'
}

foreach ($ref in $textUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$ref]
    $cell.Style = "Normal"
}
